# Slide 3, shape "Rectangle 20" (id=16) - the "Lopende acties" bullet list box
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(14)

# Append a new bullet paragraph after the existing "Hand detection" bullet.
# InsertAfter with a leading carriage return starts a new paragraph that
# inherits the formatting (bullet + run formatting) of the preceding one.
$tr = $shp.TextFrame.TextRange
$tr.InsertAfter([char]13 + " Looking for possible hardware to use for the project") | Out-Null

# The shape grew taller to fit the extra line (EMU -> points, 12700 EMU/pt).
$shp.Height = 810138 / 12700
